# Fruta / hortaliza, semanal
#
# Weekly data refresh: a new weekly price-report row is inserted into the
# "Pepino dulce" (Mapocho Venta Directa de Santiago) sheet, right before the
# row that used to be row 159. Inserting the row pushes every following row
# down by one (old row 236 becomes the new row 237), and the sheet's used
# range grows from A1:R236 to A1:R237.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 159 - this shifts rows 159..236
# down to 160..237 and extends the sheet dimension automatically.
$ws.Rows.Item(159).Insert()

# Populate the newly inserted row 159 with the new weekly record.
$ws.Cells.Item(159, 1).Value  = 12
$ws.Cells.Item(159, 2).Value  = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(159, 3).Value  = "Metropolitana"
$ws.Cells.Item(159, 4).Value  = 44777
$ws.Cells.Item(159, 5).Value  = 13
$ws.Cells.Item(159, 6).Value  = 100112043
$ws.Cells.Item(159, 7).Value  = "Pepino dulce"
$ws.Cells.Item(159, 8).Value  = "Cultivar IV Región"
$ws.Cells.Item(159, 9).Value  = "Primera"
$ws.Cells.Item(159, 10).Value = 290
$ws.Cells.Item(159, 11).Value = 15000
$ws.Cells.Item(159, 12).Value = 15000
$ws.Cells.Item(159, 13).Value = 15000
$ws.Cells.Item(159, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(159, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(159, 16).Value = 833
$ws.Cells.Item(159, 17).Value = 18
$ws.Cells.Item(159, 18).Value = "Hortaliza"
